$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Range("H17").Value = 845.1111
$ws.Range("I17").Value = 690
$ws.Range("J17").Value = 860.2439000000001
$ws.Range("K17").Value = 2070
$ws.Range("L17").Value = 2580.7317
$ws.Range("M17").Value = -1902
$ws.Range("N17").Value = -2916.7317
# row 40
$ws.Range("H40").Value = 3478.15
$ws.Range("J40").Value = 3975.7144
$ws.Range("L40").Value = 3975.7144
$ws.Range("N40").Value = -4325.7144
# row 43
$ws.Range("H43").Value = 1080.125
$ws.Range("J43").Value = 949.75
$ws.Range("L43").Value = 949.75
$ws.Range("N43").Value = -1087.75
# row 98
$ws.Range("H98").Value = 1643.6786
$ws.Range("I98").Value = 1489.4348
$ws.Range("J98").Value = 2353.2
$ws.Range("K98").Value = 1489.4348
$ws.Range("L98").Value = 2353.2
$ws.Range("M98").Value = 8.565200000000004
$ws.Range("N98").Value = -5349.2
# row 111
$ws.Range("H111").Value = 13889555
$ws.Range("I111").Value = 18519074
$ws.Range("K111").Value = 55557222
$ws.Range("M111").Value = -55554155
# row 122
$ws.Range("H122").Value = 1643.6786
$ws.Range("I122").Value = 1489.4348
$ws.Range("J122").Value = 2353.2
$ws.Range("K122").Value = 4468.3044
$ws.Range("L122").Value = 7059.599999999999
$ws.Range("M122").Value = -2018.3044
$ws.Range("N122").Value = -11959.6
# row 135
$ws.Range("H135").Value = 2999.3333
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 2999.3333
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 26993.9997
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -32063.9997

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 3415.8816
$ws.Range("I32").Value = 2762.5522
$ws.Range("K32").Value = 2762.5522
$ws.Range("M32").Value = -2475.5522
# row 55
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
# row 61
$ws.Range("H61").Value = 2635.5557
$ws.Range("I61").Value = 2496.7646
$ws.Range("K61").Value = 2496.7646
$ws.Range("M61").Value = -2284.7646
# row 101
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
# row 136
$ws.Range("H136").Value = 2635.5557
$ws.Range("I136").Value = 2496.7646
$ws.Range("K136").Value = 7490.293799999999
$ws.Range("M136").Value = -4940.293799999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 50
$ws.Range("H50").Value = 44223.168
$ws.Range("J50").Value = 44223.168
$ws.Range("L50").Value = 44223.168
$ws.Range("N50").Value = -45371.168

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 11
$ws.Range("H11").Value = 37.5
$ws.Range("I11").Value = 16.666666
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 16.666666
$ws.Range("L11").Value = 100
$ws.Range("M11").Value = 123.333334
$ws.Range("N11").Value = -380
# row 15
$ws.Range("H15").Value = 1378.1333
$ws.Range("I15").Value = 1412.2858
$ws.Range("J15").Value = 900
$ws.Range("K15").Value = 1412.2858
$ws.Range("L15").Value = 900
$ws.Range("M15").Value = -1242.2858
$ws.Range("N15").Value = -1240
# row 31
$ws.Range("H31").Value = 26303.568
$ws.Range("I31").Value = 1266.7307
$ws.Range("K31").Value = 1266.7307
$ws.Range("M31").Value = -971.7307000000001
# row 34
$ws.Range("H34").Value = 26303.568
$ws.Range("I34").Value = 1266.7307
$ws.Range("K34").Value = 1266.7307
$ws.Range("M34").Value = -1064.7307
# row 58
$ws.Range("H58").Value = 2049.9614
$ws.Range("I58").Value = 1929.3846
$ws.Range("J58").Value = 2170.5386
$ws.Range("K58").Value = 1929.3846
$ws.Range("L58").Value = 2170.5386
$ws.Range("M58").Value = -1726.3846
$ws.Range("N58").Value = -2576.5386
# row 102
$ws.Range("H102").Value = 49993.668
$ws.Range("J102").Value = 49993.668
$ws.Range("L102").Value = 49993.668
$ws.Range("N102").Value = -54861.668
# row 136
$ws.Range("H136").Value = 2049.9614
$ws.Range("I136").Value = 1929.3846
$ws.Range("J136").Value = 2170.5386
$ws.Range("K136").Value = 5788.1538
$ws.Range("L136").Value = 6511.6158
$ws.Range("M136").Value = -3238.1538
$ws.Range("N136").Value = -11611.6158

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 35763.758
$ws.Range("J5").Value = 78876.92
$ws.Range("L5").Value = 236630.76
$ws.Range("N5").Value = -236854.76
# row 22
$ws.Range("H22").Value = 213.5
$ws.Range("I22").Value = 299.5
$ws.Range("J22").Value = 127.5
$ws.Range("K22").Value = 898.5
$ws.Range("L22").Value = 382.5
$ws.Range("M22").Value = -729.5
$ws.Range("N22").Value = -720.5
# row 27
$ws.Range("H27").Value = 213.5
$ws.Range("I27").Value = 299.5
$ws.Range("J27").Value = 127.5
$ws.Range("K27").Value = 898.5
$ws.Range("L27").Value = 382.5
$ws.Range("M27").Value = -796.5
$ws.Range("N27").Value = -586.5
# row 37
$ws.Range("H37").Value = 58521.223
$ws.Range("J37").Value = 58521.223
$ws.Range("L37").Value = 175563.669
$ws.Range("N37").Value = -175787.669
# row 135
$ws.Range("H135").Value = 35763.758
$ws.Range("J135").Value = 78876.92
$ws.Range("L135").Value = 709892.28
$ws.Range("N135").Value = -714962.28
# row 138
$ws.Range("H138").Value = 3542.6365
$ws.Range("I138").Value = 2996.125
$ws.Range("K138").Value = 8988.375
$ws.Range("M138").Value = -3848.375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 18189582
$ws.Range("I70").Value = 18189582
$ws.Range("K70").Value = 18189582
$ws.Range("M70").Value = -18189312
# row 73
$ws.Range("H73").Value = 18189582
$ws.Range("I73").Value = 18189582
$ws.Range("K73").Value = 18189582
$ws.Range("M73").Value = -18188646
# row 102
$ws.Range("H102").Value = 10459345
$ws.Range("I102").Value = 27780746
$ws.Range("J102").Value = 2760945
$ws.Range("K102").Value = 27780746
$ws.Range("L102").Value = 2760945
$ws.Range("M102").Value = -27779124
$ws.Range("N102").Value = -2764189
# row 132
$ws.Range("H132").Value = 3618.611
$ws.Range("I132").Value = 2497.3845
$ws.Range("J132").Value = 6533.8
$ws.Range("K132").Value = 7492.1535
$ws.Range("L132").Value = 19601.4
$ws.Range("M132").Value = -4962.1535
$ws.Range("N132").Value = -24661.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 3230.04
$ws.Range("I7").Value = 1809.1177
$ws.Range("K7").Value = 1809.1177
$ws.Range("M7").Value = -1697.1177
# row 13
$ws.Range("H13").Value = 10000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 10000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 10000
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -10280
# row 40
$ws.Range("H40").Value = 4551.08
$ws.Range("I40").Value = 3251.5789
$ws.Range("K40").Value = 3251.5789
$ws.Range("M40").Value = -3115.5789
# row 55
$ws.Range("H55").Value = 1836.909
$ws.Range("I55").Value = 1630.7894
$ws.Range("J55").Value = 2116.6428
$ws.Range("K55").Value = 1630.7894
$ws.Range("L55").Value = 2116.6428
$ws.Range("M55").Value = -1457.7894
$ws.Range("N55").Value = -2462.6428
# row 82
$ws.Range("H82").Value = 6174066
$ws.Range("I82").Value = 7937556.5
$ws.Range("J82").Value = 1850.5
$ws.Range("K82").Value = 7937556.5
$ws.Range("L82").Value = 1850.5
$ws.Range("M82").Value = -7937195.5
$ws.Range("N82").Value = -2572.5
# row 85
$ws.Range("H85").Value = 6174066
$ws.Range("I85").Value = 7937556.5
$ws.Range("J85").Value = 1850.5
$ws.Range("K85").Value = 7937556.5
$ws.Range("L85").Value = 1850.5
$ws.Range("M85").Value = -7936308.5
$ws.Range("N85").Value = -4346.5
# row 93
$ws.Range("H93").Value = 27779956
$ws.Range("I93").Value = 41669228
$ws.Range("K93").Value = 41669228
$ws.Range("M93").Value = -41667980
# row 126
$ws.Range("H126").Value = 3230.04
$ws.Range("I126").Value = 1809.1177
$ws.Range("K126").Value = 5427.3531
$ws.Range("M126").Value = -2957.3531
# row 132
$ws.Range("H132").Value = 4476.673
$ws.Range("I132").Value = 3519.3057
$ws.Range("J132").Value = 6630.75
$ws.Range("K132").Value = 10557.9171
$ws.Range("L132").Value = 19892.25
$ws.Range("M132").Value = -8027.917099999999
$ws.Range("N132").Value = -24952.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 8
$ws.Range("H8").Value = 1262975
$ws.Range("I8").Value = 2509975
$ws.Range("J8").Value = 15975
$ws.Range("K8").Value = 2509975
$ws.Range("L8").Value = 15975
$ws.Range("M8").Value = -2509835
$ws.Range("N8").Value = -16255
# row 11
$ws.Range("H11").Value = 5009996.5
$ws.Range("I11").Value = 5009996.5
$ws.Range("K11").Value = 5009996.5
$ws.Range("M11").Value = -5009854.5
# row 13
$ws.Range("H13").Value = 3625
$ws.Range("I13").Value = 3000
$ws.Range("J13").Value = 5500
$ws.Range("K13").Value = 3000
$ws.Range("L13").Value = 5500
$ws.Range("M13").Value = -2860
$ws.Range("N13").Value = -5780
# row 81
$ws.Range("H81").Value = 11905669
$ws.Range("I81").Value = 12821336
$ws.Range("K81").Value = 25642672
$ws.Range("M81").Value = -25641611
# row 84
$ws.Range("H84").Value = 11905669
$ws.Range("I84").Value = 12821336
$ws.Range("K84").Value = 128213360
$ws.Range("M84").Value = -128208056
# row 119
$ws.Range("H119").Value = 55439.6
$ws.Range("J119").Value = 55439.6
$ws.Range("L119").Value = 55439.6
$ws.Range("N119").Value = -65115.6
